# "Generate Report for Handback"
# Update the handoff/handback timestamps recorded on the per-language
# report sheets (these values are stored as literal text, not real
# date/time serials, so assign them as strings).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 17:28:12"
$wsZhCn.Range("E3").Value = "2016-03-18 17:28:12"
$wsZhCn.Range("H2").Value = "2016-03-18 17:28:55"
$wsZhCn.Range("H3").Value = "2016-03-18 17:28:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 17:28:22"
$wsDeDe.Range("E3").Value = "2016-03-18 17:28:22"
$wsDeDe.Range("H2").Value = "2016-03-18 17:29:09"
$wsDeDe.Range("H3").Value = "2016-03-18 17:29:09"
